$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.570.77"
$ws.Range("E2").Value = "  +0.10%  "

$ws.Range("D3").Value = "2.469.12"
$ws.Range("E3").Value = "  -0.51%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.ClearFormats()
$ws.Range("E4").Value = "  +0.17%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "317.43"
$c.ClearFormats()
$ws.Range("E5").Value = "  +1.32%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "91.96"
$c.ClearFormats()
$ws.Range("E6").Value = "  -0.79%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.551"
$c.ClearFormats()
$ws.Range("E7").Value = "  +0.95%  "

$ws.Range("E8").Value = "  +0.12%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.515"
$c.ClearFormats()
$ws.Range("E9").Value = "  +0.99%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0856"
$c.ClearFormats()
$ws.Range("E10").Value = "  +8.34%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "32.85"
$c.ClearFormats()
$ws.Range("E11").Value = "  +0.11%  "

$ws.Range("E12").Value = "  +0.43%  "

$ws.Range("D13").Value = "2.852.22"
$ws.Range("E13").Value = "  -0.35%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.86"
$c.ClearFormats()
$ws.Range("E14").Value = "  -0.40%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "15.53"
$c.ClearFormats()
$ws.Range("E15").Value = "  -5.27%  "

$ws.Range("D16").Value = "2.477.14"
$ws.Range("E16").Value = "  -2.49%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.789"
$c.ClearFormats()
$ws.Range("E17").Value = "  +1.88%  "

$ws.Range("D18").Value = "41.542.62"
$ws.Range("E18").Value = "  +0.01%  "

$ws.Range("E19").Value = "  +0.07%  "

$ws.Range("E20").Value = "  -1.68%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "71.06"
$c.ClearFormats()
$ws.Range("E21").Value = "  -1.51%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "11.27"
$c.ClearFormats()
$ws.Range("E22").Value = "  +0.43%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "238.52"
$c.ClearFormats()
$ws.Range("E23").Value = "  +0.95%  "

$ws.Range("E24").Value = "  +1.00%  "

$ws.Range("E25").Value = "  +1.37%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()

$ws.Range("E27").Value = "  -0.84%  "

$ws.Range("E28").Value = "  +3.00%  "

$ws.Range("E29").Value = "  +1.46%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "35.96"
$c.ClearFormats()
$ws.Range("E30").Value = "  -0.08%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "161.82"
$c.ClearFormats()
$ws.Range("E31").Value = "  +2.42%  "

$ws.Range("E32").Value = "  +0.83%  "

$ws.Range("E33").Value = "  +0.07%  "

$ws.Range("E34").Value = "  +0.64%  "

$ws.Range("E35").Value = "  +0.80%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "17.17"
$c.ClearFormats()
$ws.Range("E36").Value = "  -2.04%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.90"
$c.ClearFormats()
$ws.Range("E37").Value = "  -0.20%  "

$ws.Range("E38").Value = "  +1.10%  "

$ws.Range("E39").Value = "  -0.18%  "

$ws.Range("E40").Value = "  -3.01%  "

$ws.Range("E41").Value = "  -2.96%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "2.47"
$c.ClearFormats()
$ws.Range("E42").Value = "  +2.87%  "

$ws.Range("D43").Value = "1.981.95"
$ws.Range("E43").Value = "  +0.51%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.0285"
$c.ClearFormats()
$ws.Range("E44").Value = "  +0.29%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "19.01"
$c.ClearFormats()
$ws.Range("E45").Value = "  -1.21%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.98"
$c.ClearFormats()
$ws.Range("E46").Value = "  +0.53%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "9.16"
$c.ClearFormats()
$ws.Range("E47").Value = "  +2.51%  "

$ws.Range("D48").Value = "2.709.75"
$ws.Range("E48").Value = "  -0.35%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "97.30"
$c.ClearFormats()
$ws.Range("E49").Value = "  -0.62%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "73.98"
$c.ClearFormats()
$ws.Range("E50").Value = "  +2.16%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "67.02"
$c.ClearFormats()
$ws.Range("E51").Value = "  -1.85%  "

